# Fix the "kafka2 was release ..." typo -> "released" and grow the
# textbox so the longer line still fits (slide 12, "KAFKA: Producer and
# Consumer"), per the target diff for ppt/slides/slide12.xml.

$p = $ppt.ActivePresentation

$target = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text.IndexOf("kafka2 was release") -ge 0) {
                $target = $shape
            }
        }
    }
}

$tr = $target.TextFrame.TextRange
$fullText = $tr.Text
$oldPhrase = "kafka2 was release in middle 2018"
$newPhrase = "kafka2 was released in middle 2018"

$startIdx = $fullText.IndexOf($oldPhrase)
$sub = $tr.Characters($startIdx + 1, $oldPhrase.Length)
$sub.Text = $newPhrase

# The line grew wider, so PowerPoint re-fit the (spAutoFit) textbox: widen
# it to keep matching the new text width.
$target.Width = 425.893937007874
